# Apply updated cryptocurrency price/volume data per commit
# "Updated cryptos list on Thu Jul 27 14:55:58 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '29.399.75'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.873.77'
$ws.Range('E3').Value = '  +0.84%  '
$r = $ws.Range('D4')
$r.NumberFormat = '@'
$r.Value = '1.002'
$r.Style = 'Normal'
$ws.Range('E4').Value = '  +0.19%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '0.7134'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  +1.54%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '241.31'
$r.Style = 'Normal'
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').Value = '  +0.16%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.07894'
$r.Style = 'Normal'
$ws.Range('E8').Value = '  -0.30%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '0.3084'
$r.Style = 'Normal'
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '25.35'
$r.Style = 'Normal'
$ws.Range('E10').Value = '  +3.74%  '
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '0.08246'
$r.Style = 'Normal'
$ws.Range('E11').Value = '  +0.81%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.7232'
$r.Style = 'Normal'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '5.249'
$r.Style = 'Normal'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.852.51'
$ws.Range('E14').Value = '  +8.75%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '91.01'
$r.Style = 'Normal'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '29.404.36'
$ws.Range('E16').Value = '  +1.84%  '
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '5.857'
$r.Style = 'Normal'
$ws.Range('E17').Value = '  +0.82%  '
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '244.49'
$r.Style = 'Normal'
$ws.Range('E18').Value = '  +2.66%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '0.000007818'
$r.Style = 'Normal'
$ws.Range('E19').Value = '  +0.05%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '13.22'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').Value = '2.116.07'
$ws.Range('E21').Value = '  +6.95%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '8.060'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  +6.88%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  +0.24%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '0.1595'
$r.Style = 'Normal'
$ws.Range('E25').Value = '  +11.31%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '162.59'
$r.Style = 'Normal'
$ws.Range('E26').Value = '  +0.31%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '9.003'
$r.Style = 'Normal'
$ws.Range('E27').Value = '  +1.19%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '18.27'
$r.Style = 'Normal'
$ws.Range('E28').Value = '  +0.83%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '1.352'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  -2.08%  '
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '1.496'
$r.Style = 'Normal'
$ws.Range('E30').Value = '  +1.33%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '4.391'
$r.Style = 'Normal'
$ws.Range('E31').Value = '  +1.58%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '4.097'
$r.Style = 'Normal'
$ws.Range('E32').Value = '  +0.83%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '0.05186'
$r.Style = 'Normal'
$ws.Range('E33').Value = '  +0.23%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '1.931'
$r.Style = 'Normal'
$ws.Range('E34').Value = '  +0.84%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '1.189'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  +1.20%  '
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '0.7211'
$r.Style = 'Normal'
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('E37').Value = '  -0.05%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '0.01859'
$r.Style = 'Normal'
$ws.Range('E38').Value = '  +0.39%  '
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '2.693'
$r.Style = 'Normal'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '1.179.85'
$ws.Range('E40').Value = '  +1.81%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.9078'
$r.Style = 'Normal'
$ws.Range('E41').Value = '  -1.63%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '6.112'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  +2.67%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '72.45'
$r.Style = 'Normal'
$ws.Range('E43').Value = '  +3.37%  '
$ws.Range('E44').Value = '  +0.15%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '101.99'
$r.Style = 'Normal'
$ws.Range('E45').Value = '  +1.02%  '
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '0.5291'
$r.Style = 'Normal'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = '2.014.48'
$ws.Range('E47').Value = '  +6.15%  '
$ws.Range('E48').Value = '  +2.07%  '
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '2.897'
$r.Style = 'Normal'
$ws.Range('E49').Value = '  +6.17%  '
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '9.262'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  +0.66%  '
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.4287'
$r.Style = 'Normal'
$ws.Range('E51').Value = '  +1.12%  '
